$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New changelog entry: date, version, changes
$ws.Range("A7").Value = 44296
$ws.Range("B7").Value = "1.0.1"
$ws.Range("C7").Value = "Two small bugfixes:`n- Selection in dropdown menu cannot be deleted`n- RAM issue, because of double entry in token data"

# Copy formatting from the previous row (row 4, which also wraps to a tall row)
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A7:C7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows.Item(7).RowHeight = 45

# Update selection to match the saved state
$ws.Range("C13:C15").Select() | Out-Null
